$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.676.91"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "3.269.60"
$ws.Range("E3").Value = "  -1.92%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.10%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").Value = "3.260.30"
$ws.Range("E9").Value = "  -2.07%  "

$ws.Range("E10").Value = "  -5.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.567"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.16%  "

$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "692.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.01%  "

$ws.Range("D15").Value = "3.790.20"
$ws.Range("E15").Value = "  -1.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.77%  "

$ws.Range("D17").Value = "66.775.90"

$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").Value = "3.269.57"
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.878"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.93%  "

$ws.Range("E24").Value = "  +2.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.75%  "

$ws.Range("E26").Value = "  -4.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "578.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.08%  "

$ws.Range("D33").Value = "3.824.63"
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.101"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.59%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.92%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -15.66%  "

$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "31.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.55%  "

$ws.Range("D43").Value = "0.0₃0658"
$ws.Range("E43").Value = "  -6.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.321"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.12%  "

$ws.Range("E46").Value = "  -4.22%  "

$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("E50").Value = "  +3.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
